# MaxsLaw.xlsx - "slight adjustment to coloring" commit
# Applies the parameter tweaks + new thermalBankFactor row to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value tweaks -------------------------------------------------

# fpsMult: 1 -> 5
$ws.Range("C6").Value = 5

# alignment: 1000 -> 1
$ws.Range("C29").Value = 1

# cohesion: 1e-6 -> 1e-3
$ws.Range("C32").Value = 0.001

# cohesionAscensionMult: 100000 -> 100
$ws.Range("C37").Value = 100

# renderScale: [150;150] -> [100;100] (row 61 pre-insert / row 62 post-insert)
$ws.Range("C61").Value = "[100;100]"

# --- Insert new "thermalBankFactor" row after collisionKillDistance (row 55) ---

$ws.Rows.Item(56).Insert()
$ws.Range("A56").Value = "Multiplier to intensify bank angle based on thermal strength"
$ws.Range("B56").Value = "thermalBankFactor"
$ws.Range("C56").Value = 0.05

# --- Visuals section tweaks (rows shifted down by 1 after the insert) ---

# showText: FALSE -> TRUE
$ws.Range("C66").Value = $true

# followAgent: FALSE -> TRUE
$ws.Range("C68").Value = $true

# followRadius: 500 -> 1500
$ws.Range("C69").Value = 1500

# --- Selection / view state ------------------------------------------------

$ws.Range("C6").Select()
